$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: activity name + hours
$ws.Range("A28").Value = "3. iterace - sekvenční diagramy a jeden scénář"
$ws.Range("A28").Style = $ws.Range("A27").Style
$ws.Range("B28").Value = 1.5

# Move selection to A29 as in the edited workbook
$ws.Range("A29").Select()
